$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New requirement rows (RNF06, RNF07, RNF08) ---
# Fill order matters: it determines the order entries land in xl/sharedStrings.xml,
# matching the target shared-string table layout.
$ws.Range("A2").Value = "RNF06"
$ws.Range("A3").Value = "RNF07"
$ws.Range("B2").Value = "El usuario puede loguearse usando servicios como Facebook, Twitter, Github"
$ws.Range("B3").Value = "Se debe informar un tiempo estimado de entrega"
$ws.Range("B4").Value = "Se pretende que se informe al usuario tanto por medio de e-mail como notificaciones push "
$ws.Range("A4").Value = "RNF08"

# B4 description ("Se pretende...") is emphasized in red font (new style, no wrap).
$ws.Range("B4").Font.Color = 255

# Row 4 grows slightly taller to fit its (now non-wrapped) red text.
$ws.Rows(4).RowHeight = 20.25

# Columns resized: narrower "Código" column, much wider "Descripción" column.
# (ColumnWidth snaps to a 1/6-character pixel grid internally, so we pick the
# input that lands as close as possible to the target stored widths.)
$ws.Columns(1).ColumnWidth = 5.666666666666667
$ws.Columns(2).ColumnWidth = 82.66666666666667

# Selection moves to C2 (just past the filled table).
$ws.Range("C2").Select()
